# Updated cryptos list on Thu Nov 28 09:08:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Some "price" strings look like plain numbers (e.g. "238.15", "1.00").
    # A straight .Value assignment lets Excel's smart-parsing coerce them to
    # numeric cells (losing the trailing zero / adding float noise), so for
    # those we briefly force Text format, assign, then restore the original
    # (default/"Normal") style so no stray number-format sticks around.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "95.423.25"
$ws.Range("E2").Value = "  +2.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.625.71"
$ws.Range("E3").Value = "  +6.23%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - Solana
Set-TextValue "D5" "238.15"
$ws.Range("E5").Value = "  +2.31%  "

# Row 6 - BNB
Set-TextValue "D6" "654.05"
$ws.Range("E6").Value = "  +5.66%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +2.86%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.405"
$ws.Range("E8").Value = "  +3.42%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.05%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.995"
$ws.Range("E10").Value = "  +2.60%  "

# Row 11 - LidoStakedEther
$ws.Range("D11").Value = "3.623.27"
$ws.Range("E11").Value = "  +6.21%  "

# Row 12 - Avalanche
Set-TextValue "D12" "42.71"
$ws.Range("E12").Value = "  -0.74%  "

# Row 13 - TRON
Set-TextValue "D13" "0.200"
$ws.Range("E13").Value = "  +0.45%  "

# Row 14 - Toncoin
$ws.Range("E14").Value = "  +1.00%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.316.00"
$ws.Range("E15").Value = "  +6.23%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "95.343.00"
$ws.Range("E16").Value = "  +2.37%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +3.19%  "

# Row 18 - was WrappedEther, now Polkadot
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "8.20"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19 - was Polkadot, now WrappedEther
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.621.19"
$ws.Range("E19").Value = "  +6.23%  "

# Row 20 - Uniswap
Set-TextValue "D20" "13.01"
$ws.Range("E20").Value = "  +11.76%  "

# Row 21 - Chainlink
Set-TextValue "D21" "17.98"
$ws.Range("E21").Value = "  -0.24%  "

# Row 22 - SuiNetwork
Set-TextValue "D22" "3.59"
$ws.Range("E22").Value = "  +5.62%  "

# Row 23 - Stellar
Set-TextValue "D23" "0.481"
$ws.Range("E23").Value = "  -1.49%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "507.56"
$ws.Range("E24").Value = "  +2.11%  "

# Row 25 - PEPE
Set-TextValue "D25" "0.0000197"
$ws.Range("E25").Value = "  +7.22%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "  -0.55%  "

# Row 27 - Litecoin
Set-TextValue "D27" "96.18"
$ws.Range("E27").Value = "  +0.87%  "

# Row 28 - Aptos
Set-TextValue "D28" "12.74"
$ws.Range("E28").Value = "  +6.22%  "

# Row 29 - WrappedeETH
$ws.Range("D29").Value = "3.814.70"
$ws.Range("E29").Value = "  +6.03%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "3.13"
$ws.Range("E30").Value = "  +14.57%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "  +0.16%  "

# Row 32 - Dai
$ws.Range("E32").Value = "  -0.06%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +0.38%  "

# Row 34 - Binance-PegBSC-USD
Set-TextValue "D34" "0.994"
$ws.Range("E34").Value = "  +1.77%  "

# Row 35 - was Cronos, now EthereumClassic
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "32.21"
$ws.Range("E35").Value = "  +11.35%  "

# Row 36 - was EthereumClassic, now Cronos
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D36" "0.177"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37 - PolygonEcosystemToken
Set-TextValue "D37" "0.562"
$ws.Range("E37").Value = "  +2.73%  "

# Row 38 - Bittensor
Set-TextValue "D38" "571.95"
$ws.Range("E38").Value = "  +2.03%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  +9.17%  "

# Row 40 - Fetch.AI
$ws.Range("E40").Value = "  +5.02%  "

# Row 41 - was USDe, now ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D41" "0.932"
$ws.Range("E41").Value = "  +4.09%  "

# Row 42 - was ARBITRUM, now USDe
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  +0.47%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "35.76"
$ws.Range("E44").Value = "  +47.92%  "

# Row 45 - ImmutableX
$ws.Range("E45").Value = "  +1.31%  "

# Row 46 - WhiteBITCoin
Set-TextValue "D46" "23.72"
$ws.Range("E46").Value = "  +0.18%  "

# Row 47 - Filecoin
$ws.Range("E47").Value = "  +4.40%  "

# Row 48 - Stacks
Set-TextValue "D48" "2.23"
$ws.Range("E48").Value = "  +5.72%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +0.24%  "

# Row 50 - MantraDAO
Set-TextValue "D50" "3.52"
$ws.Range("E50").Value = "  -5.18%  "

# Row 51 - OKB
Set-TextValue "D51" "53.78"
$ws.Range("E51").Value = "  +0.88%  "
